$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers need an explicit
# Text number format first, otherwise Excel auto-converts the string to a
# numeric value (losing trailing zeros / exact decimal text).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = "34.406.74"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "1.802.17"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "228.01"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("D6").Value = "0.584"
$ws.Range("E6").Value = "  +4.89%  "
$ws.Range("D8").Value = "35.01"
$ws.Range("E8").Value = "  +6.43%  "
$ws.Range("E9").Value = "  +0.70%  "
$ws.Range("E10").Value = "  -0.12%  "
$ws.Range("D11").Value = "0.0952"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").Value = "2.062.65"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "1.812.43"
$ws.Range("E13").Value = "  +1.05%  "
$ws.Range("D14").Value = "11.17"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Value = "0.641"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("D16").Value = "34.387.74"
$ws.Range("E16").Value = "  -0.51%  "
$ws.Range("D17").Value = "4.34"
$ws.Range("E17").Value = "  +1.28%  "
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0797"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "245.42"
$ws.Range("E20").Value = "  -1.05%  "
$ws.Range("D21").Value = "11.51"
$ws.Range("E21").Value = "  +2.24%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("D24").Value = "171.64"
$ws.Range("E24").Value = "  +3.87%  "
$ws.Range("E25").Value = "  +2.17%  "
$ws.Range("D26").Value = "7.67"
$ws.Range("E26").Value = "  +5.33%  "
$ws.Range("E27").Value = "  +2.82%  "
$ws.Range("D28").Value = "16.73"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("D30").Value = "4.00"
$ws.Range("E30").Value = "  -4.00%  "
$ws.Range("E31").Value = "  +0.93%  "
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("D35").Value = "1.394.93"
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("E36").Value = "  +0.82%  "
$ws.Range("D37").Value = "2.50"
$ws.Range("E37").Value = "  -3.23%  "
$ws.Range("D38").Value = "1.06"
$ws.Range("E38").Value = "  +0.26%  "
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("D40").Value = "83.05"
$ws.Range("E40").Value = "  -2.44%  "
$ws.Range("D41").Value = "2.83"
$ws.Range("E41").Value = "  +2.68%  "
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("D43").Value = "2.38"
$ws.Range("D44").Value = "13.51"
$ws.Range("E44").Value = "  -0.80%  "
$ws.Range("E45").Value = "  +3.07%  "
$ws.Range("D46").Value = "0.0508"
$ws.Range("E46").Value = "  -3.59%  "
$ws.Range("E47").Value = "  -2.19%  "
$ws.Range("D48").Value = "1.962.49"
$ws.Range("E48").Value = "  +0.47%  "
$ws.Range("D49").Value = "104.59"
$ws.Range("E49").Value = "  -1.40%  "
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("E51").Value = "  +0.66%  "
